$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Table 1 (SCALE A / SCALE B), rows 3-8 ---
$ws.Range("A3").Value = -47777
$ws.Range("B3").Value = 101.4
$ws.Range("F3").Value = 49996
$ws.Range("G3").Value = 101.4

$ws.Range("A4").Value = -47767
$ws.Range("B4").Value = 101.4
$ws.Range("F4").Value = 50065
$ws.Range("G4").Value = 101.4

$ws.Range("A5").Value = -47764
$ws.Range("B5").Value = 101.4
$ws.Range("F5").Value = 50048
$ws.Range("G5").Value = 101.4

$ws.Range("A6").Value = -47746
$ws.Range("B6").Value = 101.4
$ws.Range("F6").Value = 49996
$ws.Range("G6").Value = 101.4

$ws.Range("A7").Value = -47728
$ws.Range("B7").Value = 101.4
$ws.Range("F7").Value = 50046
$ws.Range("G7").Value = 101.4

$ws.Range("A8").Value = -47767
$ws.Range("B8").Value = 101.4
$ws.Range("F8").Value = 49878
$ws.Range("G8").Value = 101.4

# Row 10 used to carry a manual-offset override; it's gone now.
$ws.Range("C10").ClearContents()
$ws.Range("H10").ClearContents()

# --- Table 2 (SCALE C / SCALE D), rows 14-19 ---
$ws.Range("A14").Value = -47621
$ws.Range("B14").Value = 101.4
$ws.Range("F14").Value = -52090
$ws.Range("G14").Value = 101.4

$ws.Range("A15").Value = -47316
$ws.Range("B15").Value = 101.4
$ws.Range("F15").Value = -50224
$ws.Range("G15").Value = 101.4

$ws.Range("A16").Value = -47268
$ws.Range("B16").Value = 101.4
$ws.Range("F16").Value = -50201
$ws.Range("G16").Value = 101.4

$ws.Range("A17").Value = -47288
$ws.Range("B17").Value = 101.4
$ws.Range("F17").Value = -50189
$ws.Range("G17").Value = 101.4

$ws.Range("A18").Value = -47322
$ws.Range("B18").Value = 101.4
$ws.Range("F18").Value = -50169
$ws.Range("G18").Value = 101.4

$ws.Range("A19").Value = -47300
$ws.Range("B19").Value = 101.4
$ws.Range("F19").Value = -50166
$ws.Range("G19").Value = 101.4

# Row 21 used to carry a manual-offset override too; it's gone now.
$ws.Range("C21").ClearContents()

# Move the view / selection to reflect where the user ended up working.
$ws.Range("A20").Select()
